$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"24.35712233333334"
$ws.Cells.Item(2, 8).Value = [double]"73.07136700000001"
$ws.Cells.Item(2, 9).Value = [double]"0.3750500562097488"
$ws.Cells.Item(2, 10).Value = [double]"0.3750500562097488"
$ws.Cells.Item(2, 13).Value = [double]"15.48523133333333"
$ws.Cells.Item(2, 14).Value = [double]"46.455694"
$ws.Cells.Item(2, 15).Value = [double]"0.2291784917153818"
$ws.Cells.Item(2, 16).Value = [double]"0.2291784917153818"
$ws.Cells.Item(2, 17).Value = [double]"377.1756739459665"
$ws.Cells.Item(2, 18).Value = [double]"3394.581065513698"
$ws.Cells.Item(2, 19).Value = [double]"0.0859534061999194"
$ws.Cells.Item(2, 20).Value = [double]"0.08595340619991938"
$ws.Cells.Item(3, 7).Value = [double]"24.35712233333334"
$ws.Cells.Item(3, 8).Value = [double]"73.07136700000001"
$ws.Cells.Item(3, 9).Value = [double]"0.3750500562097488"
$ws.Cells.Item(3, 10).Value = [double]"0.3750500562097488"
$ws.Cells.Item(3, 13).Value = [double]"37.71549866666666"
$ws.Cells.Item(3, 15).Value = [double]"0.5581822391063724"
$ws.Cells.Item(3, 16).Value = [double]"0.5581822391063724"
$ws.Cells.Item(3, 17).Value = [double]"918.6410148866702"
$ws.Cells.Item(3, 18).Value = [double]"8267.769133980031"
$ws.Cells.Item(3, 19).Value = [double]"0.2093462801521284"
$ws.Cells.Item(3, 20).Value = [double]"0.2093462801521284"
$ws.Cells.Item(4, 7).Value = [double]"24.35712233333334"
$ws.Cells.Item(4, 8).Value = [double]"73.07136700000001"
$ws.Cells.Item(4, 9).Value = [double]"0.3750500562097488"
$ws.Cells.Item(4, 10).Value = [double]"0.3750500562097488"
$ws.Cells.Item(4, 13).Value = [double]"14.31939066666667"
$ws.Cells.Item(4, 14).Value = [double]"42.958172"
$ws.Cells.Item(4, 15).Value = [double]"0.2119242705923184"
$ws.Cells.Item(4, 16).Value = [double]"0.2119242705923184"
$ws.Cells.Item(4, 17).Value = [double]"348.7791502067915"
$ws.Cells.Item(4, 18).Value = [double]"3139.012351861124"
$ws.Cells.Item(4, 19).Value = [double]"0.07948220959785905"
$ws.Cells.Item(4, 20).Value = [double]"0.07948220959785905"
$ws.Cells.Item(5, 7).Value = [double]"24.35712233333334"
$ws.Cells.Item(5, 8).Value = [double]"73.07136700000001"
$ws.Cells.Item(5, 9).Value = [double]"0.3750500562097488"
$ws.Cells.Item(5, 10).Value = [double]"0.3750500562097488"
$ws.Cells.Item(5, 11).Value = [double]"1"
$ws.Cells.Item(5, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5, 13).Value = [double]"0.04831133333333334"
$ws.Cells.Item(5, 14).Value = [double]"0.144934"
$ws.Cells.Item(5, 15).Value = [double]"0.0007149985859274246"
$ws.Cells.Item(5, 16).Value = [double]"0.0007149985859274245"
$ws.Cells.Item(5, 17).Value = [double]"1.176725056086445"
$ws.Cells.Item(5, 18).Value = [double]"10.590525504778"
$ws.Cells.Item(5, 19).Value = [double]"0.0002681602598419715"
$ws.Cells.Item(5, 20).Value = [double]"0.0002681602598419715"
$ws.Cells.Item(6, 9).Value = [double]"0.2805618708302703"
$ws.Cells.Item(6, 10).Value = [double]"0.2805618708302702"
$ws.Cells.Item(6, 13).Value = [double]"15.48523133333333"
$ws.Cells.Item(6, 14).Value = [double]"46.455694"
$ws.Cells.Item(6, 15).Value = [double]"0.2291784917153818"
$ws.Cells.Item(6, 16).Value = [double]"0.2291784917153818"
$ws.Cells.Item(6, 17).Value = [double]"282.1519713485054"
$ws.Cells.Item(6, 18).Value = [double]"2539.367742136548"
$ws.Cells.Item(6, 19).Value = [double]"0.0642987463897271"
$ws.Cells.Item(6, 20).Value = [double]"0.0642987463897271"
$ws.Cells.Item(7, 9).Value = [double]"0.2805618708302703"
$ws.Cells.Item(7, 10).Value = [double]"0.2805618708302702"
$ws.Cells.Item(7, 13).Value = [double]"37.71549866666666"
$ws.Cells.Item(7, 15).Value = [double]"0.5581822391063724"
$ws.Cells.Item(7, 16).Value = [double]"0.5581822391063724"
$ws.Cells.Item(7, 17).Value = [double]"687.2033145727146"
$ws.Cells.Item(7, 18).Value = [double]"6184.829831154431"
$ws.Cells.Item(7, 19).Value = [double]"0.1566046532679131"
$ws.Cells.Item(7, 20).Value = [double]"0.1566046532679131"
$ws.Cells.Item(8, 9).Value = [double]"0.2805618708302703"
$ws.Cells.Item(8, 10).Value = [double]"0.2805618708302702"
$ws.Cells.Item(8, 13).Value = [double]"14.31939066666667"
$ws.Cells.Item(8, 14).Value = [double]"42.958172"
$ws.Cells.Item(8, 15).Value = [double]"0.2119242705923184"
$ws.Cells.Item(8, 16).Value = [double]"0.2119242705923184"
$ws.Cells.Item(8, 17).Value = [double]"260.9095219916027"
$ws.Cells.Item(8, 18).Value = [double]"2348.185697924424"
$ws.Cells.Item(8, 19).Value = [double]"0.05945786983172129"
$ws.Cells.Item(8, 20).Value = [double]"0.05945786983172129"
$ws.Cells.Item(9, 9).Value = [double]"0.2805618708302703"
$ws.Cells.Item(9, 10).Value = [double]"0.2805618708302702"
$ws.Cells.Item(9, 11).Value = [double]"1"
$ws.Cells.Item(9, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(9, 13).Value = [double]"0.04831133333333334"
$ws.Cells.Item(9, 14).Value = [double]"0.144934"
$ws.Cells.Item(9, 15).Value = [double]"0.0007149985859274246"
$ws.Cells.Item(9, 16).Value = [double]"0.0007149985859274245"
$ws.Cells.Item(9, 17).Value = [double]"0.8802669876253335"
$ws.Cells.Item(9, 18).Value = [double]"7.922402888628"
$ws.Cells.Item(9, 19).Value = [double]"0.000200601340908796"
$ws.Cells.Item(9, 20).Value = [double]"0.0002006013409087959"
$ws.Cells.Item(10, 7).Value = [double]"22.31748066666667"
$ws.Cells.Item(10, 8).Value = [double]"66.952442"
$ws.Cells.Item(10, 9).Value = [double]"0.3436437303202491"
$ws.Cells.Item(10, 10).Value = [double]"0.343643730320249"
$ws.Cells.Item(10, 13).Value = [double]"15.48523133333333"
$ws.Cells.Item(10, 14).Value = [double]"46.455694"
$ws.Cells.Item(10, 15).Value = [double]"0.2291784917153818"
$ws.Cells.Item(10, 16).Value = [double]"0.2291784917153818"
$ws.Cells.Item(10, 17).Value = [double]"345.5913509005276"
$ws.Cells.Item(10, 18).Value = [double]"3110.322158104748"
$ws.Cells.Item(10, 19).Value = [double]"0.07875575180224209"
$ws.Cells.Item(10, 20).Value = [double]"0.07875575180224208"
$ws.Cells.Item(11, 7).Value = [double]"22.31748066666667"
$ws.Cells.Item(11, 8).Value = [double]"66.952442"
$ws.Cells.Item(11, 9).Value = [double]"0.3436437303202491"
$ws.Cells.Item(11, 10).Value = [double]"0.343643730320249"
$ws.Cells.Item(11, 13).Value = [double]"37.71549866666666"
$ws.Cells.Item(11, 15).Value = [double]"0.5581822391063724"
$ws.Cells.Item(11, 16).Value = [double]"0.5581822391063724"
$ws.Cells.Item(11, 17).Value = [double]"841.7149123270257"
$ws.Cells.Item(11, 18).Value = [double]"7575.434210943232"
$ws.Cells.Item(11, 19).Value = [double]"0.191815826845023"
$ws.Cells.Item(11, 20).Value = [double]"0.191815826845023"
$ws.Cells.Item(12, 7).Value = [double]"22.31748066666667"
$ws.Cells.Item(12, 8).Value = [double]"66.952442"
$ws.Cells.Item(12, 9).Value = [double]"0.3436437303202491"
$ws.Cells.Item(12, 10).Value = [double]"0.343643730320249"
$ws.Cells.Item(12, 13).Value = [double]"14.31939066666667"
$ws.Cells.Item(12, 14).Value = [double]"42.958172"
$ws.Cells.Item(12, 15).Value = [double]"0.2119242705923184"
$ws.Cells.Item(12, 16).Value = [double]"0.2119242705923184"
$ws.Cells.Item(12, 17).Value = [double]"319.5727243617804"
$ws.Cells.Item(12, 18).Value = [double]"2876.154519256024"
$ws.Cells.Item(12, 19).Value = [double]"0.07282644689174217"
$ws.Cells.Item(12, 20).Value = [double]"0.07282644689174216"
$ws.Cells.Item(13, 7).Value = [double]"22.31748066666667"
$ws.Cells.Item(13, 8).Value = [double]"66.952442"
$ws.Cells.Item(13, 9).Value = [double]"0.3436437303202491"
$ws.Cells.Item(13, 10).Value = [double]"0.343643730320249"
$ws.Cells.Item(13, 11).Value = [double]"1"
$ws.Cells.Item(13, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(13, 13).Value = [double]"0.04831133333333334"
$ws.Cells.Item(13, 14).Value = [double]"0.144934"
$ws.Cells.Item(13, 15).Value = [double]"0.0007149985859274246"
$ws.Cells.Item(13, 16).Value = [double]"0.0007149985859274245"
$ws.Cells.Item(13, 17).Value = [double]"1.078187247647556"
$ws.Cells.Item(13, 18).Value = [double]"9.703685228828002"
$ws.Cells.Item(13, 19).Value = [double]"0.0002457047812418034"
$ws.Cells.Item(13, 20).Value = [double]"0.0002457047812418033"
$ws.Cells.Item(14, 5).Value = [double]"2"
$ws.Cells.Item(14, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(14, 7).Value = [double]"0.04834033333333334"
$ws.Cells.Item(14, 8).Value = [double]"0.145021"
$ws.Cells.Item(14, 9).Value = [double]"0.0007443426397318391"
$ws.Cells.Item(14, 10).Value = [double]"0.0007443426397318388"
$ws.Cells.Item(14, 13).Value = [double]"15.48523133333333"
$ws.Cells.Item(14, 14).Value = [double]"46.455694"
$ws.Cells.Item(14, 15).Value = [double]"0.2291784917153818"
$ws.Cells.Item(14, 16).Value = [double]"0.2291784917153818"
$ws.Cells.Item(14, 17).Value = [double]"0.7485612443971112"
$ws.Cells.Item(14, 18).Value = [double]"6.737051199574001"
$ws.Cells.Item(14, 19).Value = [double]"0.0001705873234931887"
$ws.Cells.Item(14, 20).Value = [double]"0.0001705873234931886"
$ws.Cells.Item(15, 5).Value = [double]"2"
$ws.Cells.Item(15, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(15, 7).Value = [double]"0.04834033333333334"
$ws.Cells.Item(15, 8).Value = [double]"0.145021"
$ws.Cells.Item(15, 9).Value = [double]"0.0007443426397318391"
$ws.Cells.Item(15, 10).Value = [double]"0.0007443426397318388"
$ws.Cells.Item(15, 13).Value = [double]"37.71549866666666"
$ws.Cells.Item(15, 15).Value = [double]"0.5581822391063724"
$ws.Cells.Item(15, 16).Value = [double]"0.5581822391063724"
$ws.Cells.Item(15, 17).Value = [double]"1.823179777379555"
$ws.Cells.Item(15, 18).Value = [double]"16.408617996416"
$ws.Cells.Item(15, 19).Value = [double]"0.0004154788413078658"
$ws.Cells.Item(15, 20).Value = [double]"0.0004154788413078657"
$ws.Cells.Item(16, 5).Value = [double]"2"
$ws.Cells.Item(16, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(16, 7).Value = [double]"0.04834033333333334"
$ws.Cells.Item(16, 8).Value = [double]"0.145021"
$ws.Cells.Item(16, 9).Value = [double]"0.0007443426397318391"
$ws.Cells.Item(16, 10).Value = [double]"0.0007443426397318388"
$ws.Cells.Item(16, 13).Value = [double]"14.31939066666667"
$ws.Cells.Item(16, 14).Value = [double]"42.958172"
$ws.Cells.Item(16, 15).Value = [double]"0.2119242705923184"
$ws.Cells.Item(16, 16).Value = [double]"0.2119242705923184"
$ws.Cells.Item(16, 17).Value = [double]"0.6922041179568889"
$ws.Cells.Item(16, 18).Value = [double]"6.229837061612"
$ws.Cells.Item(16, 19).Value = [double]"0.0001577442709959309"
$ws.Cells.Item(16, 20).Value = [double]"0.0001577442709959308"
$ws.Cells.Item(17, 5).Value = [double]"2"
$ws.Cells.Item(17, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(17, 7).Value = [double]"0.04834033333333334"
$ws.Cells.Item(17, 8).Value = [double]"0.145021"
$ws.Cells.Item(17, 9).Value = [double]"0.0007443426397318391"
$ws.Cells.Item(17, 10).Value = [double]"0.0007443426397318388"
$ws.Cells.Item(17, 11).Value = [double]"1"
$ws.Cells.Item(17, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(17, 13).Value = [double]"0.04831133333333334"
$ws.Cells.Item(17, 14).Value = [double]"0.144934"
$ws.Cells.Item(17, 15).Value = [double]"0.0007149985859274246"
$ws.Cells.Item(17, 16).Value = [double]"0.0007149985859274245"
$ws.Cells.Item(17, 17).Value = [double]"0.002335385957111111"
$ws.Cells.Item(17, 18).Value = [double]"0.021018473614"
$ws.Cells.Item(17, 19).Value = [double]"5.322039348537514E-07"
$ws.Cells.Item(17, 20).Value = [double]"5.322039348537511E-07"
